# Add a new "work_node" worksheet in front of the existing sheets,
# modeled on the existing "hb_node" sheet (same borders/fonts), populated
# with a single-node "offline to excel" node description + a hyperlink
# in the password cell (mirrors "add offline to excel tool").

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet by copying hb_node (keeps fonts/borders/styles
#        identical to the donor sheet instead of inventing new style ids) ---
$donor = $wb.Worksheets.Item(1)
$null = $donor.Copy($donor)
$ws = $wb.Worksheets.Item(1)
$ws.Name = "work_node"

# Drop the extra demo rows (3,4,5) that came along with the copied sheet;
# rows 1 (header) and 2 (first data row) are reused/overwritten below.
$ws.Rows.Item(3).EntireRow.Delete()
$ws.Rows.Item(3).EntireRow.Delete()
$ws.Rows.Item(3).EntireRow.Delete()

# --- 2. Header row -----------------------------------------------------
# A1 ("序号") is left untouched -- it already carries the right text/style
# from the donor sheet.
$ws.Range("B1").Value = "名称"
$ws.Range("C1").Value = "主机名"
$ws.Range("D1").Value = "IP地址"
$ws.Range("E1").Value = "数据端口"
$ws.Range("F1").Value = "缓存目录"
$ws.Range("G1").Value = "日志目录"
$ws.Range("H1").Value = "口令"
$ws.Range("I1").Value = "节点类型"

# --- 3. Data row ---------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "hadoop1"
$ws.Range("C2").Value = "hadoop1"
$ws.Range("E2").Value = 26804
$ws.Range("F2").Value = "/home/sunxo/i2data/cache/"
$ws.Range("G2").Value = "/home/sunxo/i2data/log/"
$ws.Range("H2").Value = "Info@1234"
$ws.Range("I2").Value = "源端节点|备端节点"
$ws.Range("D2").Value = "10.1.125.201"

# --- 4. Re-apply the donor-sheet borders/fonts that belong on each cell --
# (values above only change content; format-only pastes below line the
# cells back up with the rest of the table.)
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# --- 5. Hyperlink on the password cell -----------------------------------
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:Info@1234")
# Adding the hyperlink re-styles H2 with the built-in blue/underline
# "Hyperlink" style; put the plain bordered table style back.
$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# --- 6. Column widths (matches the bordered data-entry look of the rest
#        of the workbook) ---
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 11.830729166666666
$ws.Columns.Item(5).ColumnWidth = 6.330729166666667
$ws.Columns.Item(6).ColumnWidth = 21.998697916666668
$ws.Columns.Item(7).ColumnWidth = 19.830729166666668
$ws.Columns.Item(8).ColumnWidth = 8.998697916666666
$ws.Columns.Item(9).ColumnWidth = 15.998697916666666
$ws.Columns.Item(10).ColumnWidth = 5.166666666666667
$ws.Columns.Item(11).ColumnWidth = 5.330729166666667
$ws.Columns.Item(12).ColumnWidth = 7.830729166666667
$ws.Columns.Item(13).ColumnWidth = 5.166666666666667

# --- 7. View state: select I2, zoom matches the rest of the workbook ---
$ws.Range("I2").Select()

Write-Output "work_node sheet added"
